$wb = $excel.ActiveWorkbook

# Insert a new "State" column into the hotel_info sheet, between
# "Hotel_Name" and "City", and populate it with "Louisiana".
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

# Reorder the worksheet tabs so "review_info" comes before "hotel_info".
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($wb.Worksheets.Item(1))
